# Adds all IG authors as contact:
#   - Duplicate the existing "Contact" / "No display for ContactDetail" row
#     (row 11) so there are two extra Contact rows, pushing every row below
#     it down by two.
#   - Refresh the "Date" metadata value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 12..23 down to 14..25, working from the bottom up so
# that a source row is always read before it gets overwritten.
for ($r = 23; $r -ge 12; $r--) {
    $newRow = $r + 2
    # Clear the destination first: Copy() only overwrites non-blank source
    # cells, so a stale value could otherwise survive in the destination.
    $ws.Range("A" + $newRow + ":B" + $newRow).ClearContents()
    $ws.Range("A" + $r + ":B" + $r).Copy($ws.Range("A" + $newRow + ":B" + $newRow))
}

# Populate the two newly freed rows (12 and 13) with another "Contact" /
# "No display for ContactDetail" entry, copied (value + formatting) from the
# original Contact row (row 11).
$ws.Range("A11:B11").Copy($ws.Range("A12:B12"))
$ws.Range("A11:B11").Copy($ws.Range("A13:B13"))

# Update the Date property to the new commit timestamp.
$ws.Range("B8").Value2 = "2022-01-21T07:49:24+01:00"
